# Generate Report for Handback
#
# This script fills in the "handback" columns (Latest Target File, Latest
# Handback File, Latest Handback DateTime) for both locale sheets (zh-cn,
# de-de), adds hyperlinks on the newly-populated "Latest Target File" cells
# (mirroring the existing hyperlinks on column A), widens the columns that
# now hold longer text, and flips the Status column from "Ready for
# handoff" to "Handed back: in sync with en-US".

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item(1)
$zhcn     = $wb.Worksheets.Item(2)
$dede     = $wb.Worksheets.Item(3)

# ---------------------------------------------------------------------
# 1. Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
# ---------------------------------------------------------------------
$newStatus = "Handed back: in sync with en-US"

$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus

$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------
# 2. Fill in "Latest Target File" / "Latest Handback File" /
#    "Latest Handback DateTime" for zh-cn and de-de, and hyperlink the
#    newly-populated target-file cells the same way column A is linked.
# ---------------------------------------------------------------------

# zh-cn (sheet 2)
$zhcn.Range("I2").Value = "60f4349b-0f73-4053-91b3-b5a7e9dc8dee.md"
$zhcn.Range("J2").Value = "60f4349b-0f73-4053-91b3-b5a7e9dc8dee.c17e4c0197a54b8c2a9b53ea44ff7b723837bd3b.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-08-19 21:04:54"

$zhcn.Range("I3").Value = "f9e92b1e-807b-46c9-8a17-88b821d6cb7d.md"
$zhcn.Range("J3").Value = "f9e92b1e-807b-46c9-8a17-88b821d6cb7d.5a1583a7d6a4078213ead6144ca99524643607a3.zh-cn.xlf"
$zhcn.Range("K3").Value = "2016-08-19 21:04:54"

$zhcn.Hyperlinks.Add($zhcn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/126745d03ed187813b0ae6683934c83eca1b5f54/e2e/60f4349b-0f73-4053-91b3-b5a7e9dc8dee.md", "", "", "60f4349b-0f73-4053-91b3-b5a7e9dc8dee.md")
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/126745d03ed187813b0ae6683934c83eca1b5f54/e2e/f9e92b1e-807b-46c9-8a17-88b821d6cb7d.md", "", "", "f9e92b1e-807b-46c9-8a17-88b821d6cb7d.md")

# de-de (sheet 3)
$dede.Range("I2").Value = "60f4349b-0f73-4053-91b3-b5a7e9dc8dee.md"
$dede.Range("J2").Value = "60f4349b-0f73-4053-91b3-b5a7e9dc8dee.c17e4c0197a54b8c2a9b53ea44ff7b723837bd3b.de-de.xlf"
$dede.Range("K2").Value = "2016-08-19 21:05:02"

$dede.Range("I3").Value = "f9e92b1e-807b-46c9-8a17-88b821d6cb7d.md"
$dede.Range("J3").Value = "f9e92b1e-807b-46c9-8a17-88b821d6cb7d.5a1583a7d6a4078213ead6144ca99524643607a3.de-de.xlf"
$dede.Range("K3").Value = "2016-08-19 21:05:02"

$dede.Hyperlinks.Add($dede.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/126745d03ed187813b0ae6683934c83eca1b5f54/e2e/60f4349b-0f73-4053-91b3-b5a7e9dc8dee.md", "", "", "60f4349b-0f73-4053-91b3-b5a7e9dc8dee.md")
$dede.Hyperlinks.Add($dede.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/126745d03ed187813b0ae6683934c83eca1b5f54/e2e/f9e92b1e-807b-46c9-8a17-88b821d6cb7d.md", "", "", "f9e92b1e-807b-46c9-8a17-88b821d6cb7d.md")

# ---------------------------------------------------------------------
# 3. Widen the columns that now show the longer status / file-name text.
#    (COM ColumnWidth is in quantized "characters"; 29.17 / 39.17 are the
#    nearest settings that round-trip to the desired OOXML widths of
#    ~29.98 and 40 respectively.)
# ---------------------------------------------------------------------

$overview.Columns.Item(5).ColumnWidth = 29.17   # E (zh-cn status)
$overview.Columns.Item(6).ColumnWidth = 29.17   # F (de-de status)

$zhcn.Columns.Item(3).ColumnWidth  = 29.17      # C  (Status)
$zhcn.Columns.Item(9).ColumnWidth  = 39.17      # I  (Latest Target File)
$zhcn.Columns.Item(10).ColumnWidth = 39.17      # J  (Latest Handback File)

$dede.Columns.Item(3).ColumnWidth  = 29.17      # C  (Status)
$dede.Columns.Item(9).ColumnWidth  = 39.17      # I  (Latest Target File)
$dede.Columns.Item(10).ColumnWidth = 39.17      # J  (Latest Handback File)

Write-Output "Handback report generated"
